$d = $word.ActiveDocument

# The paragraph currently holds three runs whose text concatenates to
# "<id>p017r_1</id>" (split across differing run formatting: the literal
# "<id>" / "</id>" tags in Courier New, and the "p017r_1" value in Arial).
# Collapse them into a single run carrying the full string, using the
# formatting of the first ("<id>") run - this mirrors Word's own behaviour
# when a Find/Replace match spans multiple runs.
$d.Content.Find.Execute(
    "<id>p017r_1</id>",  # FindText
    $true,                # MatchCase
    $false,               # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "<id>p017r_1</id>",   # ReplaceWith
    2                     # Replace (wdReplaceAll)
) | Out-Null
